$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---

# Values that are unambiguous as text (contain two "." separators, e.g. "27.026.32")
# can be assigned directly; Excel cannot parse them as numbers so they stay text,
# matching the original inline-string cells exactly (no style change needed).
$ws.Range("D2").Value = "27.026.32"
$ws.Range("D3").Value = "1.621.65"
$ws.Range("D12").Value = "1.848.49"
$ws.Range("D13").Value = "1.620.76"
$ws.Range("D17").Value = "26.996.48"
$ws.Range("D35").Value = "1.349.52"
$ws.Range("D45").Value = "1.759.92"

# Values that look like plain decimal numbers (e.g. "148.20") would be auto-converted
# to floating point numbers by Excel (losing the trailing zero / becoming numeric).
# Force them to be stored as text by applying a text number format first, then restore
# the default "Normal" style so no extra formatting is left behind.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.251"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.540"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.116"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.768"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0178"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.848"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.98"
$ws.Range("D47").Style = "Normal"

# --- Volume(1h) (column E) updates ---
# These values always contain a leading "+"/"-" and a trailing "%", padded with two
# spaces on each side, so Excel always treats them as plain text already.
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("E16").Value = "  -3.32%  "
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("E23").Value = "  -5.36%  "
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("E28").Value = "  -2.66%  "
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("E32").Value = "  +39.46%  "
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("E35").Value = "  +3.51%  "
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("E43").Value = "  +5.33%  "
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("E46").Value = "  +31.76%  "
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("E50").Value = "  +5.87%  "
$ws.Range("E51").Value = "  +0.34%  "
